$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove row 3 (the second data row, for the "MuSCs" sending cluster entry)
$ws.Rows.Item(3).Delete()

# Update the target cluster label for row 2 from "Resolving-Mac" to "Neutrophils"
$ws.Range("D2").Value = "Neutrophils"

# Update the recalculated TPM-derived numeric values in row 2
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.05633566666666667
$ws.Range("N2").Value = 0.169007
$ws.Range("Q2").Value = 0.009472297771888889
$ws.Range("R2").Value = 0.085250679947
$ws.Range("S2").Value = 1
$ws.Range("T2").Value = 1
